$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.076.77'
$ws.Range("E2").Value = '  -0.41%  '

$ws.Range("D3").Value = '2.067.31'
$ws.Range("E3").Value = '  -0.94%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.85'
$ws.Range("E5").Value = '  +0.85%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.672'
$ws.Range("E6").Value = '  +1.52%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '57.83'
$ws.Range("E7").Value = '  +5.77%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '61.36'
$ws.Range("E9").Value = '  -0.45%  '

$ws.Range("E10").Value = '  +2.47%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0799'
$ws.Range("E11").Value = '  +7.13%  '

$ws.Range("E12").Value = '  +2.46%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '16.15'
$ws.Range("E13").Value = '  +7.42%  '

$ws.Range("D14").Value = '2.371.69'
$ws.Range("E14").Value = '  -0.84%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.807'
$ws.Range("E15").Value = '  -2.68%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.64'
$ws.Range("E16").Value = '  +9.06%  '

$ws.Range("D17").Value = '2.064.59'
$ws.Range("E17").Value = '  -1.02%  '

$ws.Range("D18").Value = '37.020.01'
$ws.Range("E18").Value = '  -0.40%  '

$ws.Range("E19").Value = '  +13.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '75.50'
$ws.Range("E20").Value = '  +3.83%  '

$ws.Range("D21").Value = '0.0₃0921'
$ws.Range("E21").Value = '  +8.81%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.45'
$ws.Range("E22").Value = '  +4.72%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '238.56'
$ws.Range("E23").Value = '  -0.67%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("E24").Value = '  -0.21%  '

$ws.Range("E25").Value = '  -2.28%  '

$ws.Range("E26").Value = '  +12.81%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '169.33'
$ws.Range("E27").Value = '  -1.68%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.32'
$ws.Range("E28").Value = '  +1.56%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.30'
$ws.Range("E29").Value = '  -1.76%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.126'
$ws.Range("E30").Value = '  +2.28%  '

$ws.Range("E31").Value = '  +5.02%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.77'
$ws.Range("E32").Value = '  +5.45%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0620'
$ws.Range("E33").Value = '  +0.23%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.49'
$ws.Range("E34").Value = '  +7.88%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0913'
$ws.Range("E35").Value = '  +1.60%  '

$ws.Range("E36").Value = '  -0.11%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.29'
$ws.Range("E37").Value = '  +2.98%  '

$ws.Range("B38").Value = 'WEMIXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.77'
$ws.Range("E38").Value = '  -1.73%  '

$ws.Range("B39").Value = 'Cronos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.116'
$ws.Range("E39").Value = '  +21.96%  '

$ws.Range("E40").Value = '  +1.33%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.92'
$ws.Range("E41").Value = '  -1.39%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0226'
$ws.Range("E42").Value = '  -0.31%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.16'
$ws.Range("E43").Value = '  -0.43%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '98.23'
$ws.Range("E44").Value = '  -0.26%  '

$ws.Range("E45").Value = '  +2.45%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.64'
$ws.Range("E46").Value = '  +14.73%  '

$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.17'
$ws.Range("E47").Value = '  -12.64%  '

$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.48'
$ws.Range("E48").Value = '  +6.26%  '

$ws.Range("D49").Value = '1.294.95'
$ws.Range("E49").Value = '  -2.14%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.91'
$ws.Range("E50").Value = '  -0.41%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.92'
$ws.Range("E51").Value = '  -0.95%  '
